$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.355.22'
$ws.Range("E2").Value = '  -0.60%  '
$ws.Range("D3").Value = '1.566.35'
$ws.Range("E3").Value = '  +0.03%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '210.96'
$ws.Range("E5").Value = '  -0.41%  '
$ws.Range("E6").Value = '  -0.61%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = '@'
$ws.Range("D8").Value = '44.43'
$ws.Range("E8").Value = '  -3.74%  '
$ws.Range("D9").NumberFormat = '@'
$ws.Range("D9").Value = '23.64'
$ws.Range("E9").Value = '  -1.71%  '
$ws.Range("E10").Value = '  -1.38%  '
$ws.Range("E11").Value = '  -0.92%  '
$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '0.0894'
$ws.Range("E12").Value = '  +0.76%  '
$ws.Range("D13").Value = '1.788.19'
$ws.Range("E13").Value = '  -0.04%  '
$ws.Range("D14").Value = '1.573.87'
$ws.Range("E14").Value = '  +0.48%  '
$ws.Range("E15").Value = '  -0.34%  '
$ws.Range("D16").Value = '28.348.68'
$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '0.514'
$ws.Range("E17").Value = '  -1.27%  '
$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '60.97'
$ws.Range("E18").Value = '  -1.97%  '
$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '228.39'
$ws.Range("E19").Value = '  -0.27%  '
$ws.Range("E20").Value = '  +0.45%  '
$ws.Range("D21").Value = '0.0₃0680'
$ws.Range("E21").Value = '  -2.03%  '
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '3.95'
$ws.Range("E23").Value = '  +1.58%  '
$ws.Range("E24").Value = '  -2.34%  '
$ws.Range("E25").Value = '  -1.76%  '
$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '150.43'
$ws.Range("E26").Value = '  -0.12%  '
$ws.Range("E28").Value = '  +0.11%  '
$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '6.35'
$ws.Range("E29").Value = '  -1.50%  '
$ws.Range("E30").Value = '  -0.03%  '
$ws.Range("E31").Value = '  +2.08%  '
$ws.Range("E32").Value = '  -3.72%  '
$ws.Range("E33").Value = '  -0.72%  '
$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '3.08'
$ws.Range("E34").Value = '  -0.42%  '
$ws.Range("D35").Value = '1.387.65'
$ws.Range("E35").Value = '  -0.26%  '
$ws.Range("E36").Value = '  +1.97%  '
$ws.Range("E37").Value = '  -3.34%  '
$ws.Range("E39").Value = '  +2.43%  '
$ws.Range("E40").Value = '  -1.85%  '
$ws.Range("E41").Value = '  -2.66%  '
$ws.Range("E42").Value = '  +2.98%  '
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '0.785'
$ws.Range("E44").Value = '  -0.44%  '
$ws.Range("E45").Value = '  -1.71%  '
$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '5.33'
$ws.Range("E46").Value = '  -3.19%  '
$ws.Range("E48").Value = '  -5.83%  '
$ws.Range("D49").Value = '1.700.90'
$ws.Range("E49").Value = '  -0.02%  '
$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '85.45'
$ws.Range("E51").Value = '  -1.92%  '
